$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = 2.24
$ws.Range("H2").Value = 3.45
$ws.Range("W2").Value = 1.75
$ws.Range("W3").Value = 1.47
$ws.Range("H4").Value = 1.99
$ws.Range("I6").Value = 10.5
$ws.Range("P6").Value = 1.51
$ws.Range("F7").Value = 3.45
$ws.Range("G7").Value = 4.8
$ws.Range("J7").Value = 2.74
$ws.Range("K7").Value = 5.4
$ws.Range("F8").Value = 1.95
$ws.Range("K8").Value = 4.2
$ws.Range("N9").Value = 2.54
$ws.Range("G13").Value = 2.26
$ws.Range("P13").Value = 1.76
$ws.Range("Q13").Value = 2.04
$ws.Range("N15").Value = 3.75
$ws.Range("W15").Value = 1.52
$ws.Range("G16").Value = 2.44
$ws.Range("J16").Value = 2.78
$ws.Range("W16").Value = 1.69
$ws.Range("Q17").Value = 1.93
$ws.Range("Q18").Value = 1.68
$ws.Range("R18").Value = 1.53
$ws.Range("F19").Value = 1.9
$ws.Range("G19").Value = 1.92
$ws.Range("I19").Value = 4.3
$ws.Range("Q19").Value = 1.67
$ws.Range("V19").Value = 1.3
$ws.Range("W19").Value = 2.08
$ws.Range("AI19").Value = 48
$ws.Range("S20").Value = 2.84
$ws.Range("F22").Value = 3.6
$ws.Range("H22").Value = 2
$ws.Range("J22").Value = 3.6
$ws.Range("T22").Value = 1.55
$ws.Range("X22").Value = 27
$ws.Range("AA22").Value = 1000
$ws.Range("AB22").Value = 21
$ws.Range("AC22").Value = 9.4
$ws.Range("AD22").Value = 11.5
$ws.Range("AH22").Value = 15
$ws.Range("J24").Value = 5.1
$ws.Range("K24").Value = 5.5
$ws.Range("F25").Value = 2.84
$ws.Range("I25").Value = 2.7
$ws.Range("G26").Value = 1.98
$ws.Range("R26").Value = 1.27
$ws.Range("W26").Value = 2.02
$ws.Range("P27").Value = 1.64
$ws.Range("Q27").Value = 2.28
$ws.Range("S27").Value = 4.3
$ws.Range("AE27").Value = 1000
$ws.Range("AK27").Value = 1000
$ws.Range("P28").Value = 1.51
$ws.Range("V29").Value = 1.71
$ws.Range("Y29").Value = 13
$ws.Range("AD29").Value = 13
$ws.Range("AJ29").Value = 75
$ws.Range("AL29").Value = 55
$ws.Range("F32").Value = 1.54
$ws.Range("Q33").Value = 1.91
$ws.Range("R35").Value = 1.51
$ws.Range("AC37").Value = 9
$ws.Range("AF37").Value = 26
$ws.Range("AJ37").Value = 70
$ws.Range("J38").Value = 2.9
$ws.Range("AD38").Value = 20
$ws.Range("F41").Value = 2.18
$ws.Range("H41").Value = 2.12
$ws.Range("Y43").Value = 1000
$ws.Range("Q44").Value = 1.33
$ws.Range("I47").Value = 3.55
$ws.Range("Q47").Value = 1.63
$ws.Range("R47").Value = 1.54
$ws.Range("S47").Value = 2.56
$ws.Range("AA47").Value = 70
$ws.Range("AN48").Value = 21
$ws.Range("W51").Value = 1.26
$ws.Range("F52").Value = 1.75
$ws.Range("V52").Value = 1.21
$ws.Range("P54").Value = 1.8
$ws.Range("Q54").Value = 1.84
$ws.Range("J56").Value = 3.65
$ws.Range("M56").Value = 1.08
$ws.Range("P56").Value = 1.84
$ws.Range("Q56").Value = 2.06
$ws.Range("T56").Value = 1.96
$ws.Range("X56").Value = 16
$ws.Range("AN56").Value = 13
$ws.Range("N57").Value = 5.2
$ws.Range("R57").Value = 1.57
$ws.Range("S57").Value = 2.68
$ws.Range("N58").Value = 3.85
$ws.Range("P60").Value = 2.24
$ws.Range("V60").Value = 1.93
$ws.Range("AG60").Value = 15.5
$ws.Range("G61").Value = 5.8
$ws.Range("I61").Value = 2.22
$ws.Range("T61").Value = 1.64
$ws.Range("U61").Value = 2
$ws.Range("V61").Value = 1.81
$ws.Range("G63").Value = 2.14
$ws.Range("H63").Value = 3.9
$ws.Range("I63").Value = 3.95
$ws.Range("O63").Value = 1.32
$ws.Range("W63").Value = 1.88
$ws.Range("Z63").Value = 27
$ws.Range("AD63").Value = 15.5
$ws.Range("AE63").Value = 46
$ws.Range("R64").Value = 1.42
$ws.Range("F60").Value = 3.9
$ws.Range("G60").Value = 3.95
$ws.Range("H60").Value = 2.04
$ws.Range("I60").Value = 2.08
$ws.Range("J60").Value = 3.85
$ws.Range("K60").Value = 3.95
